# Workbook under edit
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "cr" sheet's Regular/Augmentation columns (B, C)
# ---------------------------------------------------------------
$crSheet = $wb.Worksheets.Item("cr")

$crValues = @(
    @(0.64, 0.55000000000000004),
    @(0.64, 0.73),
    @(0.67, 0.75),
    @(0.73, 0.82),
    @(0.73, 0.81),
    @(0.79, 0.85),
    @(0.81, 0.82),
    @(0.83599999999999997, 0.83799999999999997)
)

for ($i = 0; $i -lt $crValues.Length; $i++) {
    $row = 3 + $i
    $crSheet.Cells.Item($row, 2).Value = $crValues[$i][0]
    $crSheet.Cells.Item($row, 3).Value = $crValues[$i][1]
}

# "cr" is no longer the tab-selected sheet; its stored selection becomes A2:C10
$crSheet.Range("A2:C10").Select() | Out-Null

# ---------------------------------------------------------------
# 2. Add a brand-new "Sheet1" right after "cr" with sst1/sst2 data
# ---------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $crSheet)
$newSheet.Name = "Sheet1"

$newSheet.Cells.Item(2, 1).Value = "% dataset"
$newSheet.Cells.Item(2, 2).Value = "Regular"
$newSheet.Cells.Item(2, 3).Value = "Augmentation"

$newData = @(
    @(0.001, 0.254, 0.26600000000000001),
    @(0.003, 0.214, 0.32500000000000001),
    @(0.01, 0.29499999999999998, 0.35899999999999999),
    @(0.05, 0.32600000000000001, 0.378),
    @(0.1, 0.35899999999999999, 0.39900000000000002),
    @(0.25, 0.38700000000000001, 0.41),
    @(0.5, 0.42, 0.45300000000000001),
    @(1, 0.43099999999999999, 0.434)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = 3 + $i
    $newSheet.Cells.Item($row, 1).Value = $newData[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $newData[$i][1]
    $newSheet.Cells.Item($row, 3).Value = $newData[$i][2]
}

# New sheet becomes the tab-selected / active sheet with selection E13
$newSheet.Range("E13").Select() | Out-Null
$newSheet.Activate() | Out-Null
